$wb = $excel.ActiveWorkbook

$days = $wb.Worksheets.Item("Days")
$intraday = $wb.Worksheets.Item("intraday")

# --- Widen column F on the "Days" sheet to roughly match the target stored
# width (68.54296875). The COM ColumnWidth property is internally quantized
# to steps of 1/6 character units by this runtime, so 67.66666667 is the
# closest settable value (rounds to a stored width of 68.5).
$days.Columns.Item(6).ColumnWidth = 67.66666667

# --- "Days" sheet: row 5 (Thursday / market holiday) ---
$days.Cells.Item(4, 6).Value = "the first trade should be of put side and I am observing the market "

$days.Cells.Item(5, 1).Value = "Thursday"
$days.Cells.Item(4, 2).Copy()
$days.Cells.Item(5, 2).PasteSpecial(-4122)
$days.Cells.Item(5, 2).Value = 45662
$days.Cells.Item(5, 3).Value = 0
$days.Cells.Item(5, 4).Value = 0
$days.Cells.Item(5, 6).Value = "market holiday"

# --- "Days" sheet: row 6 (Friday) ---
$days.Cells.Item(6, 1).Value = "Friday"
$days.Cells.Item(4, 2).Copy()
$days.Cells.Item(6, 2).PasteSpecial(-4122)
$days.Cells.Item(6, 2).Value = 45693
$days.Cells.Item(6, 3).Value = 832
$days.Cells.Item(6, 4).Value = 0
$days.Cells.Item(6, 6).Value = "I traded as per mysetup and it did work today as there was no hurdle"

# --- "intraday" sheet: row 162 (2025-05-02) ---
$intraday.Cells.Item(161, 1).Copy()
$intraday.Cells.Item(162, 1).PasteSpecial(-4122)
$intraday.Cells.Item(162, 1).Value = 45779
$intraday.Cells.Item(162, 2).Value = 832
$intraday.Cells.Item(162, 4).Value = 1
$intraday.Cells.Item(162, 5).Value = "this is call trading don't chase setup let it comes to you"

# --- "Days" sheet: row 8 (Monday) ---
$days.Cells.Item(8, 1).Value = "Monday"
$days.Cells.Item(4, 2).Copy()
$days.Cells.Item(8, 2).PasteSpecial(-4122)
$days.Cells.Item(8, 2).Value = 45782
$days.Cells.Item(8, 4).Value = 503
$days.Cells.Item(8, 6).Value = "I traded as per my emotion fuck again why biken why"

# --- "intraday" sheet: row 163 (2025-05-05) ---
$intraday.Cells.Item(161, 1).Copy()
$intraday.Cells.Item(163, 1).PasteSpecial(-4122)
$intraday.Cells.Item(163, 1).Value = 45782
$intraday.Cells.Item(163, 3).Value = 503
$intraday.Cells.Item(163, 4).Value = 1
$intraday.Cells.Item(163, 5).Value = "same mistake again chasing trade instead of coming trade your wau fuck man when you will control all these emoitions you can never be a trader if you do all these small mistake again and again don't regret alter saying why emotional trade so be careful as it may give pain if not traded as per pplan"

# --- Restore/update the selections shown in each sheet tab ---
[void]$days.Range("F8").Select()
[void]$intraday.Range("C163").Select()
